$wb = $excel.ActiveWorkbook
$dateFormat = "yyyy-mm-dd HH:mm:ss"

# --- Overview sheet ---
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("G2").NumberFormat = $dateFormat
$wsOverview.Range("G2").Value = "2016-10-18 13:20:32"
$wsOverview.Range("G3").NumberFormat = $dateFormat
$wsOverview.Range("G3").Value = "2016-10-18 13:20:32"

# --- zh-cn sheet ---
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").NumberFormat = $dateFormat
$wsZhCn.Range("H2").Value = "2016-10-18 13:20:20"
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("H3").Value = "2016-10-18 13:20:20"
$wsZhCn.Range("K2").NumberFormat = $dateFormat
$wsZhCn.Range("K2").Value = "2016-10-18 13:21:20"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("K3").Value = "2016-10-18 13:21:20"

# --- de-de sheet ---
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").NumberFormat = $dateFormat
$wsDeDe.Range("H2").Value = "2016-10-18 13:20:32"
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("H3").Value = "2016-10-18 13:20:32"
$wsDeDe.Range("K2").NumberFormat = $dateFormat
$wsDeDe.Range("K2").Value = "2016-10-18 13:21:38"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("K3").Value = "2016-10-18 13:21:38"
